# Landscaping Data.xlsx update
# - Fixes Growth (column H) values for rows 2-8 (were recorded as whole numbers,
#   should have been tenths).
# - Appends 7 new data rows (394-400) for 7/5/2025 (serial 45843).
# - Updates the sheet selection to U2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Correct the Growth values in H2:H8
# ---------------------------------------------------------------------------
$growthFixes = @{
    2 = 0.6
    3 = 0.2
    4 = 0.1
    5 = 0.3
    6 = 0.3
    7 = 0.2
    8 = 1.8
}
foreach ($r in $growthFixes.Keys) {
    $ws.Range("H$r").Value = $growthFixes[$r]
}

# ---------------------------------------------------------------------------
# 2. Append new rows 394-400
# ---------------------------------------------------------------------------
$newRows = @(
    @{ A=45843; B='Flowering';     C='Large';  D=67; E=89; G=0; H=0.1; I='Yes'; J=2; K='Bright';  L=8; M=0.5; N=67; O=30.1; P=7; Q=0.26; R=9.9; S=86; T=0 }
    @{ A=45843; B='Nonflowering';  C='Medium'; D=67; E=89; G=0; H=0.1; I='Yes'; J=3; K='Bright';  L=8; M=0.5; N=67; O=30.1; P=7; Q=0.26; R=9.9; S=86; T=0 }
    @{ A=45843; B='Nonflowering';  C='Small';  D=67; E=89; G=0; H=0;   I='Yes'; J=3; K='Bright';  L=8; M=0.5; N=67; O=30.1; P=7; Q=0.26; R=9.9; S=86; T=0 }
    @{ A=45843; B='Nonflowering';  C='Medium'; D=67; E=89; G=0; H=0;   I='Yes'; J=3; K='Neutral'; L=8; M=0.5; N=67; O=30.1; P=7; Q=0.26; R=9.9; S=86; T=0 }
    @{ A=45843; B='Nonflowering';  C='Medium'; D=67; E=89; G=0; H=0;   I='Yes'; J=3; K='Neutral'; L=8; M=0.5; N=67; O=30.1; P=7; Q=0.26; R=9.9; S=86; T=0 }
    @{ A=45843; B='Nonflowering';  C='Large';  D=67; E=89; G=0; H=0.2; I='Yes'; J=4; K='Neutral'; L=8; M=0.5; N=67; O=30.1; P=7; Q=0.26; R=9.9; S=86; T=0 }
    @{ A=45843; B='Tree';          C='Medium'; D=67; E=89; G=0; H=0.2; I='Yes'; J=1; K='Bright';  L=8; M=0.5; N=67; O=30.1; P=7; Q=0.26; R=9.9; S=86; T=0 }
)

$startRow = 394
$r = $startRow
foreach ($row in $newRows) {
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
    # Temp_Diff is a formula column: ABS(Low-High)
    $ws.Range("F$r").Formula = "=ABS(D$r-E$r)"
    $r = $r + 1
}
$endRow = $r - 1

# Copy the date number format from an existing date cell onto the new dates
# so the new rows reuse the same style record instead of creating a new one.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A$startRow`:A$endRow").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Update the sheet's active selection to U2 (and scroll back to the top)
# ---------------------------------------------------------------------------
$ws.Range("U2").Select() | Out-Null
